# atualizacao sabao em barra
# Inserts a new row for "SABAO EM BARRA - 200G - 200G" (code S010046) right
# before the existing "SABAO EM PO - 500G - 500G" row, which pushes every
# following row down by one (old row 79 -> new row 80, ..., old row 94 ->
# new row 95), and fixes up the ranges that depended on the old row count.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Insert a new blank row at row 79 - this shifts rows 79:94 down to 80:95
# and carries the existing cell styles (s="3"/"4") down with them, same as
# Excel's native "Insert Copied Cells"/"Insert Sheet Rows" behaviour.
$ws.Rows("79:79").Insert()

# Fill in the new row's values.
$ws.Range("A79").Value = "SABAO EM BARRA - 200G - 200G"
$ws.Range("B79").Value = "UN"
$ws.Range("C79").Value = "S010046"
$ws.Range("D79").Value = 51

# The table now spans one extra row (A1:D95 instead of A1:D94) - update the
# dependent named range and conditional formatting rule to match.
foreach ($dn in $wb.Names) {
    if ($dn.Name -eq "Planilha1!_FilterDatabase") {
        $dn.RefersTo = "=Planilha1!`$A`$1:`$D`$95"
    }
}

$fc = $ws.Range("A2:D94").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A2:D95"))

# Keep the on-screen selection roughly where the author left it (around the
# newly inserted row).
$ws.Application.ActiveWindow.ScrollRow = 76
$ws.Range("F86").Select()
